$d = $word.ActiveDocument

# Locate the "Docente(s) Responsável(eis)" heading paragraph.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Docente(s) Responsável(eis)*") {
        $target = $p
        break
    }
}

# Collapse to the end of that paragraph (right before its paragraph mark) and
# type a brand-new paragraph straight after it.
$endRange = $target.Range.Duplicate
$endRange.Collapse(0)
$endRange.InsertAfter("5701460 - Antonio Iacono`r")

# Re-acquire the freshly created paragraph (by its distinctive text) and give
# it the "ListBullet" style, matching the other instructor/requirement lists
# in this document.
$newPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*5701460 - Antonio Iacono*") {
        $newPara = $p
        break
    }
}
$newPara.Style = "ListBullet"
